$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = "NSE:5PAISA"
$ws.Range("C2").Value = "NSE:ADANIPOWER"
$ws.Range("D2").Value = "NSE:BAJAJFINSV"
$ws.Range("E2").Value = "NSE:ASIANPAINT"
$ws.Range("F2").Value = "NSE:HCLTECH"

# --- Row 3 ---
$ws.Range("B3").Value = "NSE:AXSENSEX"
$ws.Range("C3").Value = "NSE:DHARMAJ"
$ws.Range("E3").Value = "NSE:HAVELLS"
$ws.Range("F3").Value = "NSE:IDFC"

# --- Row 4 ---
$ws.Range("B4").Value = "NSE:BSLNIFTY"
$ws.Range("C4").Value = "NSE:GENSOL"

# --- Row 5 ---
$ws.Range("B5").Value = "NSE:CAMPUS"
$ws.Range("C5").Value = "NSE:HCC"

# --- Row 6 ---
$ws.Range("B6").Value = "NSE:CYBERTECH"
$ws.Range("C6").Value = "NSE:JSWHL"

# --- Row 7 ---
$ws.Range("B7").Value = "NSE:DHANI"
$ws.Range("C7").Value = "NSE:LLOYDSENGG"

# --- Row 8 ---
$ws.Range("B8").Value = "NSE:DPSCLTD"
$ws.Range("C8").Value = "NSE:MARALOVER"

# --- Row 9 ---
$ws.Range("B9").Value = "NSE:EUROTEXIND"
$ws.Range("C9").Value = "NSE:MCLEODRUSS"

# --- Row 10 ---
$ws.Range("B10").Value = "NSE:GODREJAGRO"
$ws.Range("C10").Value = "NSE:NELCO"

# --- Row 11 ---
$ws.Range("B11").Value = "NSE:GOLDBEES"
$ws.Range("C11").Value = "NSE:ORIENTALTL"

# --- Row 12 ---
$ws.Range("B12").Value = "NSE:HCLTECH"
$ws.Range("C12").Value = "NSE:PIONEEREMB"

# --- Row 13 ---
$ws.Range("B13").Value = "NSE:IDBI"
$ws.Range("C13").Value = "NSE:RAJESHEXPO"

# --- Row 14 ---
$ws.Range("B14").Value = "NSE:INDIANB"
$ws.Range("C14").Value = "NSE:REMSONSIND"

# --- Row 15 ---
$ws.Range("B15").Value = "NSE:ISMTLTD"
$ws.Range("C15").Value = "NSE:ROLEXRINGS"

# --- Row 16 ---
$ws.Range("B16").Value = "NSE:JUSTDIAL"
$ws.Range("C16").Value = "NSE:ROSSELLIND"

# --- Row 17 ---
$ws.Range("B17").Value = "NSE:KELLTONTEC"
$ws.Range("C17").Value = "NSE:RUSHIL"

# --- Row 18 ---
$ws.Range("B18").Value = "NSE:MTNL"
$ws.Range("C18").Value = "NSE:SALONA"

# --- Row 19 ---
$ws.Range("B19").Value = "NSE:NETWORK18"
$ws.Range("C19").Value = $null

# --- Row 20 ---
$ws.Range("B20").Value = "NSE:ONMOBILE"
$ws.Range("C20").Value = $null

# --- Row 21 ---
$ws.Range("B21").Value = "NSE:PGIL"
$ws.Range("C21").Value = $null

# --- Row 22 ---
$ws.Range("B22").Value = "NSE:PILITA"
$ws.Range("C22").Value = $null

# --- Row 23 ---
$ws.Range("B23").Value = "NSE:PIXTRANS"
$ws.Range("C23").Value = $null

# --- Remove now-unused rows 24-27 (data moved up into rows 2-23 above) ---
$ws.Range("A24:F27").Delete()
